$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(309).Insert()

$ws.Range("A309").Value = 8
$ws.Range("B309").Value = "Terminal La Palmera de La Serena"
$ws.Range("C309").Value = "Coquimbo"
$ws.Range("D309").Value = 44798
$ws.Range("E309").Value = 4
$ws.Range("F309").Value = 100112032
$ws.Range("G309").Value = "Zapallo italiano"
$ws.Range("H309").Value = "Bola 8"
$ws.Range("I309").Value = "Primera"
$ws.Range("J309").Value = 400
$ws.Range("K309").Value = 17000
$ws.Range("L309").Value = 18000
$ws.Range("M309").Value = 17500
$ws.Range("N309").Value = "$/caja 50 unidades"
$ws.Range("O309").Value = "Región de Arica y Parinacota"
$ws.Range("P309").Value = 350
$ws.Range("Q309").Value = 50
$ws.Range("R309").Value = "Hortaliza"
